$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the teacher/bill header details (name, designation, department)
$ws.Range("A3").Value = "নাম: Dr. K. M. Azharul Hasan"
$ws.Range("A4").Value = "পদবী: অধ্যাপক"
$ws.Range("F5").Value = "বিভাগ :সিএসই"

# Fill in the quantities used to compute individual line-item bill amounts
$ws.Range("G16").Value = 10
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1
$ws.Range("G29").Value = 15

# Fill in the total amount written out in words
$ws.Range("A32").Value = "কথায়:দশ হাজার দুইশো পাঁচ টাকা মাত্র।"

# Move the active selection to B5, matching the saved view state
$ws.Range("B5").Select()

$wb.Save()
